# Applies the numeric updates described in the commit diff for Valefor_Profits.xlsx
# (values refreshed by the scheduled market-data runner across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1022.1702
$ws.Range("I15").Value = 1022.1702
$ws.Range("K15").Value = 3066.5106
$ws.Range("M15").Value = -2897.5106

$ws.Range("H40").Value = 1533.2307
$ws.Range("I40").Value = 1544
$ws.Range("J40").Value = 1497.3334
$ws.Range("K40").Value = 1544
$ws.Range("L40").Value = 1497.3334
$ws.Range("M40").Value = -1369
$ws.Range("N40").Value = -1847.3334

$ws.Range("H111").Value = 1328.5714
$ws.Range("I111").Value = 1100
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 3300
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -233
$ws.Range("N111").Value = -10634

$ws.Range("H116").Value = 3115881.5
$ws.Range("I116").Value = 35716784
$ws.Range("J116").Value = 11033.619
$ws.Range("K116").Value = 35716784
$ws.Range("L116").Value = 11033.619
$ws.Range("M116").Value = -35713342
$ws.Range("N116").Value = -17917.619

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1850.5454
$ws.Range("I61").Value = 1418.16
$ws.Range("J61").Value = 3201.75
$ws.Range("K61").Value = 1418.16
$ws.Range("L61").Value = 3201.75
$ws.Range("M61").Value = -1206.16
$ws.Range("N61").Value = -3625.75

$ws.Range("H74").Value = 16668145
$ws.Range("I74").Value = 21740200
$ws.Range("J74").Value = 2818.2856
$ws.Range("K74").Value = 21740200
$ws.Range("L74").Value = 2818.2856
$ws.Range("M74").Value = -21739326
$ws.Range("N74").Value = -4566.2856

$ws.Range("H77").Value = 16668145
$ws.Range("I77").Value = 21740200
$ws.Range("J77").Value = 2818.2856
$ws.Range("K77").Value = 108701000
$ws.Range("L77").Value = 14091.428
$ws.Range("M77").Value = -108696632
$ws.Range("N77").Value = -22827.428

$ws.Range("H136").Value = 1850.5454
$ws.Range("I136").Value = 1418.16
$ws.Range("J136").Value = 3201.75
$ws.Range("K136").Value = 4254.48
$ws.Range("L136").Value = 9605.25
$ws.Range("M136").Value = -1704.48
$ws.Range("N136").Value = -14705.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 934.1429
$ws.Range("I107").Value = 912.05554
$ws.Range("K107").Value = 912.05554
$ws.Range("M107").Value = 1007.94446

$ws.Range("H134").Value = 6758149.5
$ws.Range("I134").Value = 8334390
$ws.Range("J134").Value = 2832.5715
$ws.Range("K134").Value = 25003170
$ws.Range("L134").Value = 8497.7145
$ws.Range("M134").Value = -25000635
$ws.Range("N134").Value = -13567.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 129400
$ws.Range("J3").Value = 6500
$ws.Range("L3").Value = 6500
$ws.Range("N3").Value = -6726

$ws.Range("H31").Value = 16132514
$ws.Range("I31").Value = 38463556
$ws.Range("J31").Value = 4539.278
$ws.Range("K31").Value = 38463556
$ws.Range("L31").Value = 4539.278
$ws.Range("M31").Value = -38463261
$ws.Range("N31").Value = -5129.278

$ws.Range("H34").Value = 16132514
$ws.Range("I34").Value = 38463556
$ws.Range("J34").Value = 4539.278
$ws.Range("K34").Value = 38463556
$ws.Range("L34").Value = 4539.278
$ws.Range("M34").Value = -38463354
$ws.Range("N34").Value = -4943.278

$ws.Range("H50").Value = 12750
$ws.Range("J50").Value = 12750
$ws.Range("L50").Value = 12750
$ws.Range("N50").Value = -14000

$ws.Range("H51").Value = 11375
$ws.Range("J51").Value = 13500
$ws.Range("L51").Value = 13500
$ws.Range("N51").Value = -14972

$ws.Range("H60").Value = 9697.667
$ws.Range("I60").Value = 4093
$ws.Range("J60").Value = 12500
$ws.Range("K60").Value = 4093
$ws.Range("L60").Value = 12500
$ws.Range("M60").Value = -3582
$ws.Range("N60").Value = -13522

$ws.Range("H61").Value = 11375
$ws.Range("J61").Value = 13500
$ws.Range("L61").Value = 13500
$ws.Range("N61").Value = -14196

$ws.Range("H134").Value = 1751.5625
$ws.Range("I134").Value = 1061.1
$ws.Range("J134").Value = 2902.3333
$ws.Range("K134").Value = 3183.3
$ws.Range("L134").Value = 8706.999899999999
$ws.Range("M134").Value = -648.2999999999997
$ws.Range("N134").Value = -13776.9999

$ws.Range("H141").Value = 36699.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 36699.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 36699.75
$ws.Range("N141").Value = -47059.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4114.706
$ws.Range("I70").Value = 4170.8
$ws.Range("J70").Value = 4034.5715
$ws.Range("K70").Value = 4170.8
$ws.Range("L70").Value = 4034.5715
$ws.Range("M70").Value = -3900.8
$ws.Range("N70").Value = -4574.5715

$ws.Range("H73").Value = 4114.706
$ws.Range("I73").Value = 4170.8
$ws.Range("J73").Value = 4034.5715
$ws.Range("K73").Value = 4170.8
$ws.Range("L73").Value = 4034.5715
$ws.Range("M73").Value = -3234.8
$ws.Range("N73").Value = -5906.5715

$ws.Range("H80").Value = 102140.5
$ws.Range("I80").Value = 2160
$ws.Range("J80").Value = 202121
$ws.Range("K80").Value = 2160
$ws.Range("L80").Value = 202121
$ws.Range("M80").Value = -1162
$ws.Range("N80").Value = -204117

$ws.Range("H83").Value = 102140.5
$ws.Range("I83").Value = 2160
$ws.Range("J83").Value = 202121
$ws.Range("K83").Value = 10800
$ws.Range("L83").Value = 1010605
$ws.Range("M83").Value = -5808
$ws.Range("N83").Value = -1020589

$ws.Range("H123").Value = 20199.2
$ws.Range("J123").Value = 20199.2
$ws.Range("L123").Value = 20199.2
$ws.Range("N123").Value = -25099.2

$ws.Range("H132").Value = 4413.923
$ws.Range("I132").Value = 5177.5835
$ws.Range("J132").Value = 3192.0667
$ws.Range("K132").Value = 15532.7505
$ws.Range("L132").Value = 9576.2001
$ws.Range("M132").Value = -13002.7505
$ws.Range("N132").Value = -14636.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2537.195
$ws.Range("I132").Value = 1814.619
$ws.Range("J132").Value = 3295.9
$ws.Range("K132").Value = 5443.857
$ws.Range("L132").Value = 9887.7
$ws.Range("M132").Value = -2913.857
$ws.Range("N132").Value = -14947.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 32944.445
$ws.Range("J82").Value = 32944.445
$ws.Range("L82").Value = 32944.445
$ws.Range("N82").Value = -33710.445

$ws.Range("H85").Value = 32944.445
$ws.Range("J85").Value = 32944.445
$ws.Range("L85").Value = 32944.445
$ws.Range("N85").Value = -35596.445

$ws.Range("H107").Value = 465
$ws.Range("I107").Value = 322.5
$ws.Range("J107").Value = 636
$ws.Range("K107").Value = 967.5
$ws.Range("L107").Value = 1908
$ws.Range("M107").Value = 952.5
$ws.Range("N107").Value = -5748

# Cell removed entirely in the source row (CRP!M141 leve profit column no longer populated)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M141").ClearContents()
